$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 - DKNY coupon
$ws.Range("A40").Value = "دكني - DKNY Coupon"
$ws.Range("C40").Value = "A21"
$ws.Range("D40").Value = "https://dkny.com.kw/"

# Row 41 - Huawei coupon
$ws.Range("A41").Value = "هواوي - Huawei Coupon"
$ws.Range("C41").Value = "AEB09"
$ws.Range("D41").Value = "https://consumer.huawei.com/ae-en/offer/"

# Row 42 - Reebok coupon
$ws.Range("A42").Value = "ريبوك - Reebok Coupon"
$ws.Range("C42").Value = "ADM84"
$ws.Range("D42").Value = "https://www.reebok.ae/"

# Images - inserted in this particular order to reproduce the shared-string table order
$ws.Range("E42").Value = "https://e.top4top.io/p_3398jj9fc1.png"
$ws.Range("E40").Value = "https://f.top4top.io/p_339868wzq2.png"
$ws.Range("E41").Value = "https://g.top4top.io/p_3398epsnm3.png"

# Remaining cells reuse already-existing shared strings
$ws.Range("B40").Value = "خصم 10% على جميع المنتجات"
$ws.Range("F40").Value = "السعودية - الإمارات - الكويت"
$ws.Range("G40").Value = "لا توجد ملاحظات"

$ws.Range("B41").Value = "خصم 10% على جميع المنتجات"
$ws.Range("F41").Value = "الإمارات"
$ws.Range("G41").Value = "لا توجد ملاحظات"

$ws.Range("B42").Value = "خصم 15% على جميع المنتجات"
$ws.Range("F42").Value = "الإمارات"
$ws.Range("G42").Value = "لا توجد ملاحظات"

# Update the view state (zoom / scroll position / selection) to match the saved workbook state
$ws.Range("D57").Select()
$excel.ActiveWindow.Zoom = 80
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
